$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row was inserted right after row 57 (duplicate of row 57's
# "Primera" Cebollín record for the Ñuble region), pushing every subsequent
# row down by one (old row 58 -> new row 59, ..., old row 72 -> new row 73).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the same record as row 57.
$ws.Range("A58").Value = 11
$ws.Range("B58").Value = "Vega Monumental Concepción"
$ws.Range("C58").Value = "Bíobío"
$ws.Range("D58").Value = 44358
$ws.Range("E58").Value = 8
$ws.Range("F58").Value = 100112037
$ws.Range("G58").Value = "Cebollín"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 200
$ws.Range("K58").Value = 600
$ws.Range("L58").Value = 700
$ws.Range("M58").Value = 650
$ws.Range("N58").Value = "$/paquete 6 unidades"
$ws.Range("O58").Value = "Región de Ñuble"
$ws.Range("P58").Value = 108
$ws.Range("Q58").Value = 6
$ws.Range("R58").Value = "Hortaliza"
